$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # quality_comparison
$ws2 = $wb.Worksheets.Item(2)   # computational_comparison

$xlPasteFormats = -4122

# --- Build the two new border styles once, on sheet 1, then replicate them
# by copy/paste-format onto every other cell that needs the same look.
# (Re-deriving the style from scratch on every cell would leave throw-away
# style entries behind in the workbook's style table.)

# C1: interior cell of the merged header row -> keep only top+bottom border
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.LineStyle = 1              # full box (matches the existing boxed style)
$c1.Borders.Item(7).LineStyle = -4142  # drop left
$c1.Borders.Item(10).LineStyle = -4142 # drop right

# D1: right-most cell of the merged header row -> keep top+bottom+right border
$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.LineStyle = 1              # full box
$d1.Borders.Item(7).LineStyle = -4142  # drop left

# Replicate onto sheet 2's matching header cells (two merged header blocks)
$c1.Copy() | Out-Null
$ws2.Range("C1").PasteSpecial($xlPasteFormats) | Out-Null
$ws2.Range("F1").PasteSpecial($xlPasteFormats) | Out-Null

$d1.Copy() | Out-Null
$ws2.Range("D1").PasteSpecial($xlPasteFormats) | Out-Null
$ws2.Range("G1").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# --- Anonymize "fedcore" -> "approach" ---
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Drop the stray empty inline-string cell ---
$ws2.Range("G5").ClearContents()
